$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Update the summary figures at the top of the statement
# ---------------------------------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 1912900
# Cant. Periodos (count of periods covered)
$ws.Range("F13").Value = 20

# ---------------------------------------------------------------------------
# 2. Insert a new data row right after the existing period "2406" row for
#    ALEXIS (old row 21), pushing the remaining data rows (and the
#    signature block further down) down by one row. This new blank row
#    will host MARCELA's "2406" record, and the former last detail row
#    (which carries the special "closing" bottom-border formatting) now
#    sits one row lower and will be re-used for a brand-new "2508" record.
# ---------------------------------------------------------------------------
$ws.Rows("22:22").Insert()

# Copy the formatting of a normal data row into the freshly inserted row
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Write out the full, final contents of the worker detail table
#    (rows 16 to 36). Period codes and document numbers are text values,
#    names/types are text as well; only "Valor Mora" and "Salario Basico"
#    are numeric. Row 22 is reserved for MARCELA's record (set below), all
#    other rows (16-21, 23-36) hold ALEXIS's 20 monthly periods in order.
# ---------------------------------------------------------------------------
$periods = @("2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")
$pIndex = 0
for ($r = 16; $r -le 36; $r++) {
    if ($r -eq 22) { continue }
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1102122349"
    $ws.Cells.Item($r, 4).Value = "ALEXIS JOSE ARCIA CHICA"
    $ws.Cells.Item($r, 5).Value = $periods[$pIndex]
    $ws.Cells.Item($r, 6).Value = 83600
    $ws.Cells.Item($r, 7).Value = 2090000
    $pIndex++
}

# The new MARCELA row (row 22) with its own mora/salario figures
$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "53039048"
$ws.Cells.Item(22, 4).Value = "MARCELA ANDREA CHAMORRO LEON"
$ws.Cells.Item(22, 5).Value = "2406"
$ws.Cells.Item(22, 6).Value = 240900
$ws.Cells.Item(22, 7).Value = 6335670

# ---------------------------------------------------------------------------
# 4. Refresh the worksheet dimension/selection bookkeeping so downstream
#    consumers see the correct used range.
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
